$d = $word.ActiveDocument

# Helper: replace the text of a single paragraph in-place via Range.InsertXML,
# rebuilding the paragraph's exact original XML (pPr / leading empty run /
# rPr all preserved) with just the <w:t> content swapped. This avoids the
# run-merging that a plain Find/Replace triggers when the target run has no
# run-level formatting and sits next to an empty <w:r/>.
function Set-ParagraphXml($paragraph, $innerXml) {
    $wrapped = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($wrapped) | Out-Null
}

# Finds the (first) paragraph whose text equals $oldText and rewrites it to
# $newXmlText (already the desired inner <w:p>...</w:p> markup).
function Replace-BulletParagraph($oldText, $newText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($oldText + "`r")) {
            $innerXml = '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p>'
            Set-ParagraphXml $p $innerXml
            return
        }
    }
}

# 1. Title heading + its repeated bold echo near the end of the document.
#    wdReplaceAll (2) replaces every occurrence of the exact phrase in one
#    call, which covers both the Heading1 paragraph and the bold paragraph
#    near the end (neither has a leading empty run, so plain Find/Replace
#    is structurally safe here).
$d.Content.Find.Execute(
    "Play Lucky 88 for Free: A Review of Aristocrat's Slot Game",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Lucky 88 for Free - Exciting Asian Slot Game", 2) | Out-Null

# 2. "What we like" bullet list
Replace-BulletParagraph "Variety of bonuses and symbols" "Numerous symbols and bonuses"
Replace-BulletParagraph "Bright, celebratory visual design" "Bright and lively visuals"
Replace-BulletParagraph "Culturally significant concept of luck" "Incorporates Asian cultural themes"
Replace-BulletParagraph "Lucrative bonus features" "Exciting bonus features"

# 3. "What we don't like" bullet list (content swapped between the two bullets)
Replace-BulletParagraph "Limited paylines compared to other games" "Low-value symbols"
Replace-BulletParagraph "Some may find the cultural theme too specific" "Limited number of paylines"

# 4. Meta description paragraph (italic) at the very end of the document
$d.Content.Find.Execute(
    "Read our review of Aristocrat's Lucky 88 slot game, play for free, and learn about its gameplay mechanics, visuals, bonus features, and cultural significance.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Lucky 88, an Asian-themed slot game with exciting bonus features. Play for free!", 2) | Out-Null
